$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (prices / 1h volume % changes, and two reordered coin rows).
$ws.Range("D2").Value = '52.233.40'
$ws.Range("E2").Value = '  -0.42%  '
$ws.Range("D3").Value = '2.935.95'
$ws.Range("E3").Value = '  +0.64%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '356.23'
$ws.Range("E5").Value = '  +0.46%  '
$ws.Range("D6").NumberFormat = "@"  # preserve exact text (avoid trailing-zero strip)
$ws.Range("D6").Value = '109.80'
$ws.Range("E6").Value = '  -2.39%  '
$ws.Range("E7").Value = '  +1.66%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").Value = '0.626'
$ws.Range("E9").Value = '  -0.45%  '
$ws.Range("D10").Value = '38.98'
$ws.Range("E10").Value = '  -2.88%  '
$ws.Range("E11").Value = '  +1.54%  '
$ws.Range("D12").Value = '0.0868'
$ws.Range("E12").Value = '  +0.37%  '
$ws.Range("D13").Value = '19.51'
$ws.Range("E13").Value = '  -1.77%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '7.79'
$ws.Range("E14").Value = '  -0.07%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '3.401.65'
$ws.Range("E15").Value = '  +0.72%  '
$ws.Range("D16").Value = '2.944.04'
$ws.Range("E16").Value = '  +1.23%  '
$ws.Range("D17").Value = '0.978'
$ws.Range("E17").Value = '  -2.59%  '
$ws.Range("D18").Value = '52.180.39'
$ws.Range("E18").Value = '  -0.49%  '
$ws.Range("D19").Value = '3.55'
$ws.Range("E19").Value = '  +7.78%  '
$ws.Range("D20").Value = '7.57'
$ws.Range("E20").Value = '  -1.42%  '
$ws.Range("D21").Value = '13.86'
$ws.Range("E21").Value = '  -2.65%  '
$ws.Range("E22").Value = '  -0.25%  '
$ws.Range("D23").Value = '70.39'
$ws.Range("E23").Value = '  -0.81%  '
$ws.Range("D24").Value = '270.19'
$ws.Range("E24").Value = '  +0.31%  '
$ws.Range("D25").Value = '2.79'
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("D26").Value = '0.178'
$ws.Range("E26").Value = '  +2.98%  '
$ws.Range("D27").Value = '7.88'
$ws.Range("E27").Value = '  +18.94%  '
$ws.Range("E28").Value = '  +0.26%  '
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("D30").Value = '0.107'
$ws.Range("E30").Value = '  +7.80%  '
$ws.Range("E31").Value = '  -1.85%  '
$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").Value = '37.56'
$ws.Range("E32").Value = '  -1.37%  '
$ws.Range("B33").Value = 'Toncoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D33").Value = '2.28'
$ws.Range("E33").Value = '  +0.65%  '
$ws.Range("D34").NumberFormat = "@"  # preserve exact text (avoid trailing-zero strip)
$ws.Range("D34").Value = '6.20'
$ws.Range("E34").Value = '  -1.47%  '
$ws.Range("D35").Value = '52.01'
$ws.Range("E35").Value = '  -2.51%  '
$ws.Range("D36").Value = '0.0442'
$ws.Range("E36").Value = '  -2.44%  '
$ws.Range("E37").Value = '  +0.13%  '
$ws.Range("D38").Value = '3.18'
$ws.Range("E38").Value = '  -5.76%  '
$ws.Range("D39").Value = '18.16'
$ws.Range("E39").Value = '  -4.41%  '
$ws.Range("E40").Value = '  -3.96%  '
$ws.Range("D41").Value = '2.73'
$ws.Range("E41").Value = '  -2.54%  '
$ws.Range("E42").Value = '  +1.94%  '
$ws.Range("D43").Value = '22.88'
$ws.Range("E43").Value = '  -1.25%  '
$ws.Range("D44").Value = '119.83'
$ws.Range("E44").Value = '  -0.90%  '
$ws.Range("E45").Value = '  -1.05%  '
$ws.Range("D46").Value = '3.45'
$ws.Range("E46").Value = '  -2.27%  '
$ws.Range("E47").Value = '  -5.33%  '
$ws.Range("D48").Value = '2.133.84'
$ws.Range("E48").Value = '  -2.61%  '
$ws.Range("D49").Value = '0.249'
$ws.Range("E49").Value = '  -6.13%  '
$ws.Range("E50").Value = '  +1.27%  '
$ws.Range("B51").Value = 'SEI'
$ws.Range("C51").Value = 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'
$ws.Range("D51").Value = '0.918'
$ws.Range("E51").Value = '  -4.39%  '
